# Automatische test-sync: 2025-08-28 21:18:50
# Appends the new "Klacht levering" log entry to the Logs sheet and the
# matching "Klacht / Probleem" aggregate row to the Dashboard sheet, then
# widens the chart's category/value series references and the existing
# conditional-formatting rules so they keep covering the full data range.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: new row 27 -------------------------------------------------
$logs.Range("A27").Value = "Klacht levering"
$logs.Range("B27").Value = "mailmind.test@zohomail.eu"
$logs.Range("D27").Value = "Klacht / Probleem"
$logs.Range("F27").Value = "2025-08-28 21:18:34"
$logs.Range("G27").Value = "Nee"
$logs.Range("H27").Value = "Ja"
$logs.Range("I27").Value = "Nee"
$logs.Range("J27").Value = "Nee"

# --- Extend existing conditional formatting ranges to include row 27 ------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`26")
    $newRange = $logs.Range("$col`2:$col`27")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: new row 7 ---------------------------------------------
$dashboard.Range("A7").Value = "Klacht / Probleem"
$dashboard.Range("B7").Value = 1

# --- Update chart series references to cover the extra row -----------------
$chart = $dashboard.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
